# Added Slim the Knight
# Fill in the combat stats for "Slim" (row 5 of the Characters sheet) and
# move the active selection to match the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Characters")

# Health, MP, Attack, Defense, Resistance, Skill, Speed for Slim the Knight
$ws.Range("G5").Value2 = 85
$ws.Range("H5").Value2 = 40
$ws.Range("I5").Value2 = 70
$ws.Range("J5").Value2 = 85
$ws.Range("K5").Value2 = 15
$ws.Range("L5").Value2 = 55
$ws.Range("M5").Value2 = 35

# Recalculate so the Total column (shared formula in N5) picks up the new stats
$excel.Calculate()

# Reflect where the editor ended up after entering the new row of data
$ws.Range("L6").Select() | Out-Null
